$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.422.13'
$ws.Range("E2").Value = '  +2.56%  '
$ws.Range("D3").Value = '2.065.04'
$ws.Range("E3").Value = '  +3.94%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.30'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +6.83%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.384'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.32'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("E11").Value = '  +2.28%  '
$ws.Range("E12").Value = '  +3.37%  '
$ws.Range("D13").Value = '2.364.31'
$ws.Range("E13").Value = '  +3.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.04'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.779'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.31%  '
$ws.Range("D18").Value = '2.060.11'
$ws.Range("E18").Value = '  +3.54%  '
$ws.Range("D19").Value = '37.519.47'
$ws.Range("E19").Value = '  +2.91%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.23'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +18.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.09'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("E22").Value = '  +1.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '226.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.33%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("E26").Value = '  +1.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("E28").Value = '  +13.83%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.10%  '
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("E32").Value = '  +2.59%  '
$ws.Range("E33").Value = '  +3.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0625'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.09%  '
$ws.Range("E35").Value = '  +9.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.49'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.75%  '
$ws.Range("E37").Value = '  +5.27%  '
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0985'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.67%  '
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("E43").Value = '  +26.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.05'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +9.67%  '
$ws.Range("D45").Value = '1.473.84'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("E46").Value = '  +6.21%  '
$ws.Range("E47").Value = '  +4.47%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.28%  '
$ws.Range("E49").Value = '  +3.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.28'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.56%  '
$ws.Range("E51").Value = '  +1.61%  '
